# Generate Report for Handoff
# Adds two new handed-off files (two .png images) to the localization status
# report, refreshes the already-recorded .md handoff's timestamp, and fixes
# up every dependent sheet (Overview, zh-cn, de-de) + hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: (re)point a cell's hyperlink, replacing whatever was there before
# so we don't end up with stacked/duplicate hyperlink objects on one cell.
# ---------------------------------------------------------------------------
function Set-CellHyperlink {
    param(
        $Sheet,
        [string]$CellRef,
        [string]$Url,
        [string]$Display
    )
    $range = $Sheet.Range($CellRef)
    $range.Hyperlinks.Delete()
    $Sheet.Hyperlinks.Add($range, $Url, "", "", $Display) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 2 (existing .md handoff) -> becomes the .png handoff, date refreshed
$wsOverview.Range("A2").Value = "7d0065ed-c0d5-4a16-a90d-47f36af942b0.png"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-43-18 12:43:57"
Set-CellHyperlink -Sheet $wsOverview -CellRef "A2" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/7d0065ed-c0d5-4a16-a90d-47f36af942b0.png" `
    -Display "7d0065ed-c0d5-4a16-a90d-47f36af942b0.png"

# Row 3 (new second .png handoff)
$wsOverview.Range("A3").Value = "d6eb20d4-8248-4dc0-a9e9-24db5f611ffc.png"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-43-18 12:43:57"
Set-CellHyperlink -Sheet $wsOverview -CellRef "A3" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/d6eb20d4-8248-4dc0-a9e9-24db5f611ffc.png" `
    -Display "d6eb20d4-8248-4dc0-a9e9-24db5f611ffc.png"

# Row 4 (the original .md file, now listed in its own row)
$wsOverview.Range("A4").Value = "fbffd695-54fb-415f-b6af-74689e380d87.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-43-18 12:43:57"
Set-CellHyperlink -Sheet $wsOverview -CellRef "A4" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/fbffd695-54fb-415f-b6af-74689e380d87.md" `
    -Display "fbffd695-54fb-415f-b6af-74689e380d87.md"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 2 -> first .png dependency row
$wsZhCn.Range("A2").Value = "7d0065ed-c0d5-4a16-a90d-47f36af942b0.png"
$wsZhCn.Range("B2").Value = ".png"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("D2").Value = "89192824c4ede300a47bb7d6749d7be2b322f09a.png"
$wsZhCn.Range("E2").Value = "2016-03-18 12:43:55"
$wsZhCn.Range("F2").Value = ""
$wsZhCn.Range("G2").Value = ""
$wsZhCn.Range("H2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I2").Value = "IsDependency"
$wsZhCn.Range("J2").Value = "e2e\fbffd695-54fb-415f-b6af-74689e380d87.md"
$wsZhCn.Range("K2").Value = ""
Set-CellHyperlink -Sheet $wsZhCn -CellRef "A2" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/7d0065ed-c0d5-4a16-a90d-47f36af942b0.png" `
    -Display "7d0065ed-c0d5-4a16-a90d-47f36af942b0.png"
Set-CellHyperlink -Sheet $wsZhCn -CellRef "B2" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/7d0065ed-c0d5-4a16-a90d-47f36af942b0.png" `
    -Display ".png"
Set-CellHyperlink -Sheet $wsZhCn -CellRef "D2" `
    -Url "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b7dd966120ee67869bd89288aaaba81afc87cbb5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/89192824c4ede300a47bb7d6749d7be2b322f09a.png" `
    -Display "89192824c4ede300a47bb7d6749d7be2b322f09a.png"

# Row 3 -> second .png dependency row
$wsZhCn.Range("A3").Value = "d6eb20d4-8248-4dc0-a9e9-24db5f611ffc.png"
$wsZhCn.Range("B3").Value = ".png"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "fcbea161834932b0c728f9df17eae24fd2c5008e.png"
$wsZhCn.Range("E3").Value = "2016-03-18 12:43:55"
$wsZhCn.Range("F3").Value = ""
$wsZhCn.Range("G3").Value = ""
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I3").Value = "IsDependency"
$wsZhCn.Range("J3").Value = "e2e\fbffd695-54fb-415f-b6af-74689e380d87.md"
$wsZhCn.Range("K3").Value = ""
Set-CellHyperlink -Sheet $wsZhCn -CellRef "A3" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/d6eb20d4-8248-4dc0-a9e9-24db5f611ffc.png" `
    -Display "d6eb20d4-8248-4dc0-a9e9-24db5f611ffc.png"
Set-CellHyperlink -Sheet $wsZhCn -CellRef "B3" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/d6eb20d4-8248-4dc0-a9e9-24db5f611ffc.png" `
    -Display ".png"
Set-CellHyperlink -Sheet $wsZhCn -CellRef "D3" `
    -Url "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b7dd966120ee67869bd89288aaaba81afc87cbb5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/fcbea161834932b0c728f9df17eae24fd2c5008e.png" `
    -Display "fcbea161834932b0c728f9df17eae24fd2c5008e.png"

# Row 4 -> the .md file (not a dependency; Include)
$wsZhCn.Range("A4").Value = "fbffd695-54fb-415f-b6af-74689e380d87.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "fbffd695-54fb-415f-b6af-74689e380d87.7a008ebf22e072ccd1a73a38756fe812347c4aa1.zh-cn.xlf"
$wsZhCn.Range("E4").Value = "2016-03-18 12:43:55"
$wsZhCn.Range("F4").Value = ""
$wsZhCn.Range("G4").Value = ""
$wsZhCn.Range("H4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I4").Value = "Include"
$wsZhCn.Range("J4").Value = ""
$wsZhCn.Range("K4").Value = ""
Set-CellHyperlink -Sheet $wsZhCn -CellRef "A4" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/fbffd695-54fb-415f-b6af-74689e380d87.md" `
    -Display "fbffd695-54fb-415f-b6af-74689e380d87.md"
Set-CellHyperlink -Sheet $wsZhCn -CellRef "B4" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/fbffd695-54fb-415f-b6af-74689e380d87.md" `
    -Display ".md"
Set-CellHyperlink -Sheet $wsZhCn -CellRef "D4" `
    -Url "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b7dd966120ee67869bd89288aaaba81afc87cbb5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/fbffd695-54fb-415f-b6af-74689e380d87.7a008ebf22e072ccd1a73a38756fe812347c4aa1.zh-cn.xlf" `
    -Display "fbffd695-54fb-415f-b6af-74689e380d87.7a008ebf22e072ccd1a73a38756fe812347c4aa1.zh-cn.xlf"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 2 -> first .png dependency row
$wsDeDe.Range("A2").Value = "7d0065ed-c0d5-4a16-a90d-47f36af942b0.png"
$wsDeDe.Range("B2").Value = ".png"
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("D2").Value = "89192824c4ede300a47bb7d6749d7be2b322f09a.png"
$wsDeDe.Range("E2").Value = "2016-03-18 12:43:57"
$wsDeDe.Range("F2").Value = ""
$wsDeDe.Range("G2").Value = ""
$wsDeDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I2").Value = "IsDependency"
$wsDeDe.Range("J2").Value = "e2e\fbffd695-54fb-415f-b6af-74689e380d87.md"
$wsDeDe.Range("K2").Value = ""
Set-CellHyperlink -Sheet $wsDeDe -CellRef "A2" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/7d0065ed-c0d5-4a16-a90d-47f36af942b0.png" `
    -Display "7d0065ed-c0d5-4a16-a90d-47f36af942b0.png"
Set-CellHyperlink -Sheet $wsDeDe -CellRef "B2" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/7d0065ed-c0d5-4a16-a90d-47f36af942b0.png" `
    -Display ".png"
Set-CellHyperlink -Sheet $wsDeDe -CellRef "D2" `
    -Url "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e861d851e5fffb2759ebcf41050aaa0261d1822/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/89192824c4ede300a47bb7d6749d7be2b322f09a.png" `
    -Display "89192824c4ede300a47bb7d6749d7be2b322f09a.png"

# Row 3 -> second .png dependency row
$wsDeDe.Range("A3").Value = "d6eb20d4-8248-4dc0-a9e9-24db5f611ffc.png"
$wsDeDe.Range("B3").Value = ".png"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "fcbea161834932b0c728f9df17eae24fd2c5008e.png"
$wsDeDe.Range("E3").Value = "2016-03-18 12:43:57"
$wsDeDe.Range("F3").Value = ""
$wsDeDe.Range("G3").Value = ""
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I3").Value = "IsDependency"
$wsDeDe.Range("J3").Value = "e2e\fbffd695-54fb-415f-b6af-74689e380d87.md"
$wsDeDe.Range("K3").Value = ""
Set-CellHyperlink -Sheet $wsDeDe -CellRef "A3" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/d6eb20d4-8248-4dc0-a9e9-24db5f611ffc.png" `
    -Display "d6eb20d4-8248-4dc0-a9e9-24db5f611ffc.png"
Set-CellHyperlink -Sheet $wsDeDe -CellRef "B3" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/d6eb20d4-8248-4dc0-a9e9-24db5f611ffc.png" `
    -Display ".png"
Set-CellHyperlink -Sheet $wsDeDe -CellRef "D3" `
    -Url "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e861d851e5fffb2759ebcf41050aaa0261d1822/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/fcbea161834932b0c728f9df17eae24fd2c5008e.png" `
    -Display "fcbea161834932b0c728f9df17eae24fd2c5008e.png"

# Row 4 -> the .md file (not a dependency; Include)
$wsDeDe.Range("A4").Value = "fbffd695-54fb-415f-b6af-74689e380d87.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "fbffd695-54fb-415f-b6af-74689e380d87.7a008ebf22e072ccd1a73a38756fe812347c4aa1.de-de.xlf"
$wsDeDe.Range("E4").Value = "2016-03-18 12:43:57"
$wsDeDe.Range("F4").Value = ""
$wsDeDe.Range("G4").Value = ""
$wsDeDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I4").Value = "Include"
$wsDeDe.Range("J4").Value = ""
$wsDeDe.Range("K4").Value = ""
Set-CellHyperlink -Sheet $wsDeDe -CellRef "A4" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/fbffd695-54fb-415f-b6af-74689e380d87.md" `
    -Display "fbffd695-54fb-415f-b6af-74689e380d87.md"
Set-CellHyperlink -Sheet $wsDeDe -CellRef "B4" `
    -Url "https://github.com/OpenLocalizationTest/oltest/blob/9e276cfa1b86e09d20b35303abf9ee3ad5c986bd/e2e/fbffd695-54fb-415f-b6af-74689e380d87.md" `
    -Display ".md"
Set-CellHyperlink -Sheet $wsDeDe -CellRef "D4" `
    -Url "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e861d851e5fffb2759ebcf41050aaa0261d1822/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/fbffd695-54fb-415f-b6af-74689e380d87.7a008ebf22e072ccd1a73a38756fe812347c4aa1.de-de.xlf" `
    -Display "fbffd695-54fb-415f-b6af-74689e380d87.7a008ebf22e072ccd1a73a38756fe812347c4aa1.de-de.xlf"

# NOTE: keep a trivial trailing statement after the last Hyperlinks.Add call.
# (The COM-interop host tries to surface the final top-level statement's
# result; when that happens to be a raw Hyperlink COM object coming out of
# Hyperlinks.Add, it throws "You cannot call a method on a null-valued
# expression." This no-op sidesteps it without changing any written data.)
Write-Host "Report regenerated: Overview/zh-cn/de-de updated."
